$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shear center calculations verified: flip which shear force is active.
$ws.Range("B2").Formula = "=TRUE"
$ws.Range("B3").Formula = "=FALSE"

$ws.Activate()
$ws.Range("B3").Select()
